# A new weekly price record was inserted as row 21 of the data table
# (Fruta / Vega Monumental Concepción - Mango). Every existing record
# from row 21 down to row 130 shifts down by one row (to rows 22-131),
# and the worksheet dimension grows from A1:T130 to A1:T131.
#
# The brand-new record re-uses the Volumen/Precio/Unidad/Precio-$kg
# values that used to sit in row 21 (i.e. columns other than the date
# and the origin), but carries a new Fecha (D) and a new Origen (R).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 21:130 down to 22:131, leaving row 21 blank.
$ws.Rows.Item(21).Insert()

# Row 22 now holds what used to be row 21's data; duplicate it into the
# freshly emptied row 21 so every column besides Fecha/Origen keeps its
# prior value.
$ws.Range("A22:T22").Copy()
$ws.Range("A21").PasteSpecial(-4104)
$ws.Application.CutCopyMode = $false

# Overwrite the new record's Fecha (D) and Origen (R).
$ws.Range("D21").Value2 = 44819
$ws.Range("R21").Value2 = "Brasil"
